# Update Hardik Pandya's match-by-match batting activity (runs/balls/fours/sixes)
# for Mumbai Indians. Values are stored as text (numbers-as-text), so each new
# value is written with a leading apostrophe to force a text cell instead of
# turning it into a numeric cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = "'8"
$ws.Range("D2").Value = "'4"
$ws.Range("F2").Value = "'1"

# Row 3
$ws.Range("C3").Value = "'21"
$ws.Range("F3").Value = "'1"

# Row 4
$ws.Range("C4").Value = "'17"
$ws.Range("D4").Value = "'15"
$ws.Range("E4").Value = "'0"
$ws.Range("F4").Value = "'2"

# Row 5
$ws.Range("C5").Value = "'60"
$ws.Range("D5").Value = "'21"
$ws.Range("E5").Value = "'2"
$ws.Range("F5").Value = "'7"

# Row 6
$ws.Range("C6").Value = "'37"
$ws.Range("D6").Value = "'14"
$ws.Range("F6").Value = "'5"

# Row 7
$ws.Range("C7").Value = "'30"
$ws.Range("D7").Value = "'19"
$ws.Range("E7").Value = "'2"

# Row 9
$ws.Range("C9").Value = "'14"
$ws.Range("D9").Value = "'10"
$ws.Range("E9").Value = "'0"
$ws.Range("F9").Value = "'2"

# Row 11
$ws.Range("C11").Value = "'18"
$ws.Range("D11").Value = "'13"
$ws.Range("F11").Value = "'1"

# Row 12
$ws.Range("C12").Value = "'28"
$ws.Range("D12").Value = "'19"
$ws.Range("E12").Value = "'2"
$ws.Range("F12").Value = "'2"

# Row 13
$ws.Range("C13").Value = "'15"
$ws.Range("D13").Value = "'13"
$ws.Range("F13").Value = "'1"

# Row 14
$ws.Range("D14").Value = "'11"
$ws.Range("E14").Value = "'3"
$ws.Range("F14").Value = "'2"
